$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cryptocurrency price/volume data (and the Maker/HuobiToken row swap)
$updates = @{
    "D2" = "26.463.98"
    "E2" = "  -0.78%  "
    "D3" = "1.627.03"
    "E3" = "  -0.67%  "
    "E4" = "  +0.20%  "
    "D5" = "213.04"
    "E5" = "  -0.01%  "
    "D6" = "0.499"
    "E6" = "  +1.52%  "
    "E7" = "  +0.23%  "
    "E9" = "  -1.40%  "
    "D10" = "18.82"
    "E10" = "  -1.27%  "
    "D11" = "0.0845"
    "E11" = "  +1.07%  "
    "D12" = "1.853.13"
    "E12" = "  -0.71%  "
    "D13" = "1.642.66"
    "E13" = "  +0.36%  "
    "D14" = "4.13"
    "E14" = "  +1.91%  "
    "D15" = "0.522"
    "E15" = "  -0.28%  "
    "D16" = "64.78"
    "E16" = "  +2.76%  "
    "D17" = "26.511.00"
    "D18" = "0.0₃0740"
    "E18" = "  -0.11%  "
    "D19" = "214.76"
    "E19" = "  +2.23%  "
    "E20" = "  +0.25%  "
    "E21" = "  -0.34%  "
    "E22" = "  +1.82%  "
    "D23" = "9.28"
    "E23" = "  -1.25%  "
    "D24" = "2.01"
    "E24" = "  +3.35%  "
    "D25" = "148.71"
    "E25" = "  +1.84%  "
    "E26" = "  +0.25%  "
    "E27" = "  -0.47%  "
    "D28" = "6.83"
    "E28" = "  +1.75%  "
    "E29" = "  +0.87%  "
    "E30" = "  -1.58%  "
    "E31" = "  -0.90%  "
    "E32" = "  +2.91%  "
    "D33" = "2.95"
    "E33" = "  -0.45%  "
    "E34" = "  -0.45%  "
    "B35" = "Maker"
    "C35" = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
    "D35" = "1.218.07"
    "E35" = "  +4.11%  "
    "B36" = "HuobiToken"
    "C36" = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
    "D36" = "2.38"
    "E36" = "  -0.96%  "
    "E37" = "  +3.70%  "
    "E38" = "  +0.24%  "
    "D39" = "0.794"
    "E39" = "  -1.55%  "
    "D40" = "0.506"
    "E40" = "  +0.64%  "
    "E41" = "  -2.37%  "
    "E42" = "  -0.54%  "
    "E43" = "  -0.23%  "
    "D44" = "1.763.46"
    "E44" = "  -0.52%  "
    "D45" = "93.06"
    "E45" = "  +0.81%  "
    "D46" = "1.59"
    "E46" = "  +1.95%  "
    "E47" = "  +0.13%  "
    "E48" = "  -0.74%  "
    "E49" = "  -0.56%  "
    "D50" = "7.51"
    "E50" = "  -0.65%  "
    "E51" = "  -0.75%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
